$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q2").Value = 486058
$ws.Range("R2").Value = 6546491

$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
